$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Il10-Il10ra rows (ECs/FAPs/M2/sCs sending clusters x ECs/M2 target clusters),
# following updated computations (Dr Hou advice).
$data = @(
    @("ECs", "Il10", "Il10ra", "ECs", 2, 0.6666666666666666, 36.26745966666667, 108.802379, 0.643144256662053, 0.6431442566620529, 1, 0.3333333333333333, 17.50728866666667, 52.521866, 0.2429365863908956, 0.2429365863908957, 634.9448855910239, 5714.503970319214, 0.1562432702703892, 0.1562432702703892),
    @("ECs", "Il10", "Il10ra", "M2", 2, 0.6666666666666666, 36.26745966666667, 108.802379, 0.643144256662053, 0.6431442566620529, 3, 1, 54.557973, 163.673919, 0.7570634136091043, 0.7570634136091045, 1978.679085272589, 17808.1117674533, 0.4869009863916638, 0.4869009863916638),
    @("FAPs", "Il10", "Il10ra", "ECs", 3, 1, 3.209141666666667, 9.627425, 0.05690889438359308, 0.05690889438359308, 1, 0.3333333333333333, 17.50728866666667, 52.521866, 0.2429365863908956, 0.2429365863908957, 56.18336953056112, 505.6503257750501, 0.01382525253683012, 0.01382525253683012),
    @("FAPs", "Il10", "Il10ra", "M2", 3, 1, 3.209141666666667, 9.627425, 0.05690889438359308, 0.05690889438359308, 3, 1, 54.557973, 163.673919, 0.7570634136091043, 0.7570634136091045, 175.084264403175, 1575.758379628575, 0.04308364184676296, 0.04308364184676296),
    @("M2", "Il10", "Il10ra", "ECs", 3, 1, 15.98639866666667, 47.95919600000001, 0.2834927116945642, 0.2834927116945642, 1, 0.3333333333333333, 17.50728866666667, 52.521866, 0.2429365863908956, 0.2429365863908957, 279.8784961977485, 2518.906465779736, 0.06887075164577576, 0.06887075164577577),
    @("M2", "Il10", "Il10ra", "M2", 3, 1, 15.98639866666667, 47.95919600000001, 0.2834927116945642, 0.2834927116945642, 3, 1, 54.557973, 163.673919, 0.7570634136091043, 0.7570634136091045, 872.1855068232361, 7849.669561409126, 0.2146219600487884, 0.2146219600487884),
    @("sCs", "Il10", "Il10ra", "ECs", 3, 1, 0.927863, 2.783589, 0.01645413725978976, 0.01645413725978976, 1, 0.3333333333333333, 17.50728866666667, 52.521866, 0.2429365863908956, 0.2429365863908957, 16.24436538411933, 146.199288457074, 0.00399731193790057, 0.003997311937900571),
    @("sCs", "Il10", "Il10ra", "M2", 3, 1, 0.927863, 2.783589, 0.01645413725978976, 0.01645413725978976, 3, 1, 54.557973, 163.673919, 0.7570634136091043, 0.7570634136091045, 50.622324501699, 455.600920515291, 0.01245682532188919, 0.01245682532188919)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}
